$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Bmp6"
$ws.Range("C2").Value = "Acvr2a"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 4.145393333333334
$ws.Range("H2").Value = 12.43618
$ws.Range("I2").Value = 0.1621900462138432
$ws.Range("J2").Value = 0.1621900462138432
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 7.940402333333334
$ws.Range("N2").Value = 23.821207
$ws.Range("O2").Value = 0.1931648990487216
$ws.Range("P2").Value = 0.1931648990487216
$ws.Range("Q2").Value = 32.91609089658445
$ws.Range("R2").Value = 296.24481806926
$ws.Range("S2").Value = 0.03132942390360453
$ws.Range("T2").Value = 0.03132942390360452

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Bmp6"
$ws.Range("C3").Value = "Acvr2a"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 4.145393333333334
$ws.Range("H3").Value = 12.43618
$ws.Range("I3").Value = 0.1621900462138432
$ws.Range("J3").Value = 0.1621900462138432
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 26.95384733333333
$ws.Range("N3").Value = 80.861542
$ws.Range("O3").Value = 0.6557019380820612
$ws.Range("P3").Value = 0.6557019380820612
$ws.Range("Q3").Value = 111.7342990432845
$ws.Range("R3").Value = 1005.60869138956
$ws.Range("S3").Value = 0.1063483276400361
$ws.Range("T3").Value = 0.1063483276400361

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Bmp6"
$ws.Range("C4").Value = "Acvr2a"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 4.145393333333334
$ws.Range("H4").Value = 12.43618
$ws.Range("I4").Value = 0.1621900462138432
$ws.Range("J4").Value = 0.1621900462138432
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 6.212609666666666
$ws.Range("N4").Value = 18.637829
$ws.Range("O4").Value = 0.1511331628692172
$ws.Range("P4").Value = 0.1511331628692172
$ws.Range("Q4").Value = 25.75371069480222
$ws.Range("R4").Value = 231.78339625322
$ws.Range("S4").Value = 0.02451229467020263
$ws.Range("T4").Value = 0.02451229467020263

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Bmp6"
$ws.Range("C5").Value = "Acvr2a"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 6.772999666666666
$ws.Range("H5").Value = 20.318999
$ws.Range("I5").Value = 0.2649961151116367
$ws.Range("J5").Value = 0.2649961151116367
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 7.940402333333334
$ws.Range("N5").Value = 23.821207
$ws.Range("O5").Value = 0.1931648990487216
$ws.Range("P5").Value = 0.1931648990487216
$ws.Range("Q5").Value = 53.78034235686589
$ws.Range("R5").Value = 484.023081211793
$ws.Range("S5").Value = 0.05118794782384272
$ws.Range("T5").Value = 0.05118794782384272

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Bmp6"
$ws.Range("C6").Value = "Acvr2a"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 6.772999666666666
$ws.Range("H6").Value = 20.318999
$ws.Range("I6").Value = 0.2649961151116367
$ws.Range("J6").Value = 0.2649961151116367
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 26.95384733333333
$ws.Range("N6").Value = 80.861542
$ws.Range("O6").Value = 0.6557019380820612
$ws.Range("P6").Value = 0.6557019380820612
$ws.Range("Q6").Value = 182.5583990040509
$ws.Range("R6").Value = 1643.025591036458
$ws.Range("S6").Value = 0.1737584662629172
$ws.Range("T6").Value = 0.1737584662629172

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Bmp6"
$ws.Range("C7").Value = "Acvr2a"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 6.772999666666666
$ws.Range("H7").Value = 20.318999
$ws.Range("I7").Value = 0.2649961151116367
$ws.Range("J7").Value = 0.2649961151116367
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 6.212609666666666
$ws.Range("N7").Value = 18.637829
$ws.Range("O7").Value = 0.1511331628692172
$ws.Range("P7").Value = 0.1511331628692172
$ws.Range("Q7").Value = 42.07800320146344
$ws.Range("R7").Value = 378.702028813171
$ws.Range("S7").Value = 0.0400497010248768
$ws.Range("T7").Value = 0.04004970102487681

# Row 8
$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Bmp6"
$ws.Range("C8").Value = "Acvr2a"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 14.640471
$ws.Range("H8").Value = 43.921413
$ws.Range("I8").Value = 0.5728138386745201
$ws.Range("J8").Value = 0.5728138386745202
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 7.940402333333334
$ws.Range("N8").Value = 23.821207
$ws.Range("O8").Value = 0.1931648990487216
$ws.Range("P8").Value = 0.1931648990487216
$ws.Range("Q8").Value = 116.251230089499
$ws.Range("R8").Value = 1046.261070805491
$ws.Range("S8").Value = 0.1106475273212744
$ws.Range("T8").Value = 0.1106475273212744

# Row 9
$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Bmp6"
$ws.Range("C9").Value = "Acvr2a"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 14.640471
$ws.Range("H9").Value = 43.921413
$ws.Range("I9").Value = 0.5728138386745201
$ws.Range("J9").Value = 0.5728138386745202
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 26.95384733333333
$ws.Range("N9").Value = 80.861542
$ws.Range("O9").Value = 0.6557019380820612
$ws.Range("P9").Value = 0.6557019380820612
$ws.Range("Q9").Value = 394.617020222094
$ws.Range("R9").Value = 3551.553181998846
$ws.Range("S9").Value = 0.3755951441791079
$ws.Range("T9").Value = 0.3755951441791081

# Row 10
$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Bmp6"
$ws.Range("C10").Value = "Acvr2a"
$ws.Range("D10").Value = "sCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 14.640471
$ws.Range("H10").Value = 43.921413
$ws.Range("I10").Value = 0.5728138386745201
$ws.Range("J10").Value = 0.5728138386745202
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 6.212609666666666
$ws.Range("N10").Value = 18.637829
$ws.Range("O10").Value = 0.1511331628692172
$ws.Range("P10").Value = 0.1511331628692172
$ws.Range("Q10").Value = 90.955531659153
$ws.Range("R10").Value = 818.599784932377
$ws.Range("S10").Value = 0.08657116717413772
$ws.Range("T10").Value = 0.08657116717413775
